# Adds a new "2022-Q3" quarterly sheet to the workbook (right after "总计"),
# fills it with fund-holding data, and updates the "总计" (totals) summary
# sheet with a new row for 2022-Q3 while shifting the existing rows down.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Update the "总计" (totals) sheet: insert a new first data row for
#    2022-Q3 and push the existing quarters down by one row.
# ---------------------------------------------------------------------
$totals = $wb.Worksheets.Item(1)

# Give the new bottom row (row 9) the same look (border/bold) as the
# existing column-A cells before we populate it, by copying formatting
# from the row directly above (same-sheet copy preserves style+value).
$totals.Cells.Item(8, 1).Copy($totals.Cells.Item(9, 1))

# Shift existing data rows (2..8) down to (3..9), working from the
# bottom up so that we never overwrite a row before reading it.
for ($r = 8; $r -ge 2; $r--) {
    $dst = $r + 1
    $totals.Cells.Item($dst, 2).Value2 = $totals.Cells.Item($r, 2).Value2
    $totals.Cells.Item($dst, 3).Value2 = $totals.Cells.Item($r, 3).Value2
    $totals.Cells.Item($dst, 4).Value2 = $totals.Cells.Item($r, 4).Value2
}

# Fill in the new row for 2022-Q3.
$totals.Cells.Item(2, 2).Value2 = "2022-Q3"
$totals.Cells.Item(2, 3).Value2 = 7
$totals.Cells.Item(2, 4).Value2 = 0.85

# Re-sequence column A (row index) values: 0..7 for rows 2..9.
for ($r = 2; $r -le 9; $r++) {
    $totals.Cells.Item($r, 1).Value2 = $r - 2
}

# ---------------------------------------------------------------------
# 2. Insert the new "2022-Q3" worksheet right after "总计".
# ---------------------------------------------------------------------
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $totals)
$newSheet.Name = "2022-Q3"

# Header row.
$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $col = $i + 2   # headers start at column B
    $cell = $newSheet.Cells.Item(1, $col)
    $cell.NumberFormat = "@"
    $cell.Value2 = $headers[$i]
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108   # xlCenter
    $cell.VerticalAlignment = -4160     # xlTop
    $cell.Borders.LineStyle = 1
}

# Fund holdings data for 2022-Q3.
$rows = @(
    @("483003", "工银精选平衡混合",               "15.78", "65.64", "2.50", "0.3945", 7),
    @("000893", "工银创新动力股票",                "11.92", "81.96", "2.98", "0.3552", 9),
    @("160135", "南方中证高铁产业指数（LOF）",      "1.84",  "95.01", "2.66", "0.0489", 6),
    @("160639", "鹏华中证高铁产业指数（LOF）A",     "0.75",  "94.62", "2.63", "0.0197", 6),
    @("008629", "大成景瑞稳健配置混合A",            "0.89",  "21.44", "1.59", "0.0142", 5),
    @("008630", "大成景瑞稳健配置混合C",            "0.75",  "21.44", "1.59", "0.0119", 5),
    @("015678", "鹏华中证高铁产业指数（LOF）C",     "0.06",  "94.62", "2.63", "0.0016", 6)
)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $i + 2
    $data = $rows[$i]

    $aCell = $newSheet.Cells.Item($r, 1)
    $aCell.Value2 = $i
    $aCell.Font.Bold = $true
    $aCell.HorizontalAlignment = -4108
    $aCell.VerticalAlignment = -4160
    $aCell.Borders.LineStyle = 1

    for ($c = 0; $c -lt 6; $c++) {
        $cell = $newSheet.Cells.Item($r, $c + 2)   # B..G
        $cell.NumberFormat = "@"
        $cell.Value2 = $data[$c]
    }

    $newSheet.Cells.Item($r, 8).Value2 = $data[6]  # H column, numeric
}
